$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 15; existing rows 15-118 shift down to 16-119
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new data record
$ws.Cells.Item(15, 1).Value = 11
$ws.Cells.Item(15, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(15, 3).Value = "Bíobío"
$ws.Cells.Item(15, 4).Value = 44819
$ws.Cells.Item(15, 5).Value = 8
$ws.Cells.Item(15, 6).Value = 100112001
$ws.Cells.Item(15, 7).Value = "Berenjena"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 220
$ws.Cells.Item(15, 11).Value = 12000
$ws.Cells.Item(15, 12).Value = 13000
$ws.Cells.Item(15, 13).Value = 12455
$ws.Cells.Item(15, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(15, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(15, 16).Value = 208
$ws.Cells.Item(15, 17).Value = 60
$ws.Cells.Item(15, 18).Value = "Hortaliza"
